# Merge the "zircon" sheet data into the "apatite" sheet (already present in
# column B), rename the remaining sheet to "size", update the shared-string
# header labels to be more descriptive, delete the now-redundant "zircon"
# sheet, and move the active-cell selection.

$wb = $excel.ActiveWorkbook

$apatiteSheet = $wb.Worksheets.Item("apatite")
$zirconSheet = $wb.Worksheets.Item("zircon")

# Update header labels in the remaining sheet to the more descriptive names.
# (Update B1 before A1 so the shared-string table keeps its original slot
# order: zircon label stays index 0, apatite label stays index 1.)
$apatiteSheet.Range("B1").Value = "zircon_obs_w_max"
$apatiteSheet.Range("A1").Value = "apatite_obs_w_max"

# Remove the redundant zircon sheet (its single column of data already
# exists in column B of the apatite sheet).
$null = $zirconSheet.Delete()

# Rename the remaining sheet to reflect the consolidated data.
$apatiteSheet.Name = "size"

# Move the selection as recorded in the saved view state.
$null = $apatiteSheet.Range("C18").Select()

$null = $wb.Save()
